$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.636.57"
$ws.Range("E2").Value = "  -2.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.292.36"
$ws.Range("E3").Value = "  -5.00%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.65"
$ws.Range("E5").Value = "  -1.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.96"
$ws.Range("E6").Value = "  -4.05%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  -2.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.290.00"
$ws.Range("E9").Value = "  -5.07%  "

$ws.Range("E10").Value = "  -3.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -2.88%  "

$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.332"
$ws.Range("E13").Value = "  -5.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.90"
$ws.Range("E14").Value = "  -3.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.705.02"
$ws.Range("E15").Value = "  -4.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.619.75"
$ws.Range("E16").Value = "  -2.24%  "

$ws.Range("E17").Value = "  -3.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.267.04"
$ws.Range("E18").Value = "  -6.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.63"
$ws.Range("E19").Value = "  -5.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -5.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.93"
$ws.Range("E21").Value = "  -3.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.42"
$ws.Range("E22").Value = "  -5.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.21"
$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  -6.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  -6.31%  "

$ws.Range("E28").Value = "  -6.78%  "

$ws.Range("E29").Value = "  -1.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.05"
$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0725"
$ws.Range("E31").Value = "  -6.05%  "

$ws.Range("E32").Value = "  +1.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.77"
$ws.Range("E33").Value = "  -5.62%  "

$ws.Range("E34").Value = "  -5.27%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.69"
$ws.Range("E36").Value = "  -3.91%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("E38").Value = "  -5.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.96"
$ws.Range("E39").Value = "  -6.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.07"
$ws.Range("E40").Value = "  -2.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.50"
$ws.Range("E41").Value = "  -5.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "299.98"
$ws.Range("E42").Value = "  -7.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "140.23"
$ws.Range("E43").Value = "  -4.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0948"
$ws.Range("E45").Value = "  -1.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0499"
$ws.Range("E46").Value = "  -3.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.553"
$ws.Range("E47").Value = "  -4.32%  "

$ws.Range("E48").Value = "  -6.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0214"
$ws.Range("E49").Value = "  -3.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.56"
$ws.Range("E50").Value = "  -4.66%  "

$ws.Range("E51").Value = "  -0.36%  "
